$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A74").Value = "GRT-USD"
